# Weekly update: insert a new price record (week of 2023-02-03) as the
# newest row for this product/market, pushing the existing historical
# rows (72-103) down by one (to 73-104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 72; this shifts rows 72:103 -> 73:104
# and extends the sheet dimension to A1:R104 automatically.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record.
$ws.Cells.Item(72, 1).Value = 1
$ws.Cells.Item(72, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value = 44960
$ws.Cells.Item(72, 5).Value = 15
$ws.Cells.Item(72, 6).Value = 100112040
$ws.Cells.Item(72, 7).Value = "Cilantro"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 450
$ws.Cells.Item(72, 11).Value = 1700
$ws.Cells.Item(72, 12).Value = 2000
$ws.Cells.Item(72, 13).Value = 1867
$ws.Cells.Item(72, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 934
$ws.Cells.Item(72, 17).Value = 2
$ws.Cells.Item(72, 18).Value = "Hortaliza"
